$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.527.81"
$ws.Range("E2").Value = "  +2.24%  "
$ws.Range("D3").Value = "'1.858.93"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("D5").Value = "'245.73"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "'0.6960"
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("D7").Value = "'1.0000"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.3080"
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").Value = "'0.07702"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "'23.68"
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("D11").Value = "'0.07786"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "'5.161"
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("D13").Value = "'1.851.61"
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("D14").Value = "'0.6962"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").Value = "'91.35"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "'6.336"
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("D17").Value = "'29.516.35"
$ws.Range("E17").Value = "  +2.21%  "
$ws.Range("D18").Value = "'0.000008313"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "'2.099.60"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").Value = "'238.23"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "'12.77"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'7.649"
$ws.Range("E23").Value = "  +2.49%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("D26").Value = "'8.912"
$ws.Range("E26").Value = "  +1.29%  "
$ws.Range("D27").Value = "'160.02"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").Value = "'18.30"
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").Value = "'1.538"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").Value = "'4.255"
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("D31").Value = "'4.158"
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").Value = "'1.209"
$ws.Range("E32").Value = "  +2.46%  "
$ws.Range("D33").Value = "'0.05119"
$ws.Range("E33").Value = "  +0.33%  "
$ws.Range("D34").Value = "'0.7792"
$ws.Range("E34").Value = "  +2.34%  "
$ws.Range("D35").Value = "'1.886"
$ws.Range("E35").Value = "  +2.58%  "
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("D37").Value = "'2.686"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").Value = "'1.314.09"
$ws.Range("E38").Value = "  +7.68%  "
$ws.Range("D39").Value = "'0.01876"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("D40").Value = "'2.726"
$ws.Range("E40").Value = "  +1.11%  "
$ws.Range("D41").Value = "'0.9484"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").Value = "'105.81"
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("D43").Value = "'5.771"
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "'9.830"
$ws.Range("E45").Value = "  +3.17%  "
$ws.Range("E46").Value = "  +2.04%  "
$ws.Range("D47").Value = "'1.993.76"
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("D48").Value = "'0.5229"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("D50").Value = "'63.11"
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("D51").Value = "'6.970"
$ws.Range("E51").Value = "  +1.07%  "
